{"js": "// The recorded change collapses the three blank spacer paragraphs that sat\n// between the main text blocks of this \"Ausblick\" brick:\n//   - the blank paragraph between \"... Abschied nehmen.\" and\n//     \"Freudentage, Regenmomente ...\"\n//   - the blank paragraph between \"... hei\u00dft es: BIBELVERS\" and\n//     \"Gott verhei\u00dft uns seinen Schutz ...\"\n//   - the blank paragraph between \"... des Trostes und des Lebens.\" and\n//     \"Ich finde dieses Motiv ...\"\n// Every paragraph in the document shares the same (Arial) run/paragraph\n// formatting, so removing each blank paragraph simply lets the following\n// text paragraph take its place - reproducing the same final text flow as\n// the diff (which re-homed each content run one paragraph earlier) without\n// having to rebuild any runs by hand.\n\nconst body = context.document.body;\n\nasync function deleteBlankParagraphBefore(anchorText) {\n  const results = body.search(anchorText, { matchCase: false, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    return;\n  }\n\n  const anchorParagraph = results.items[0].paragraphs.getFirst();\n  const previousParagraph = anchorParagraph.getPrevious();\n  previousParagraph.load(\"text\");\n  await context.sync();\n\n  if (previousParagraph.text === \"\") {\n    previousParagraph.delete();\n    await context.sync();\n  }\n}\n\nawait deleteBlankParagraphBefore(\"Freudentage, Regenmomente und alles dazwischen\");\nawait deleteBlankParagraphBefore(\"Gott verhei\u00dft uns seinen Schutz\");\nawait deleteBlankParagraphBefore(\"Ich finde dieses Motiv des Regenbogens\");\n", "ps1": "# The recorded change collapses the three blank spacer paragraphs that sat\n# between the main text blocks of this \"Ausblick\" brick:\n#   - the blank paragraph between \"... Abschied nehmen.\" and\n#     \"Freudentage, Regenmomente ...\"\n#   - the blank paragraph between \"... hei\u00dft es: BIBELVERS\" and\n#     \"Gott verhei\u00dft uns seinen Schutz ...\"\n#   - the blank paragraph between \"... des Trostes und des Lebens.\" and\n#     \"Ich finde dieses Motiv ...\"\n# Every paragraph in the document shares the same (Arial) run/paragraph\n# formatting, so removing each blank paragraph simply lets the following\n# text paragraph take its place - reproducing the same final text flow as\n# the diff (which re-homed each content run one paragraph earlier) without\n# having to rebuild any runs by hand.\n\n$d = $word.ActiveDocument\n$wdParagraph = 4\n\nfunction Remove-BlankParagraphBefore {\n    param([string]$AnchorText)\n\n    $searchRange = $d.Content\n    $found = $searchRange.Find.Execute($AnchorText)\n    if (-not $found) {\n        return\n    }\n\n    # Expand the (collapsed, find-sized) range to its full paragraph so we\n    # reliably know where that paragraph starts, then look just before it.\n    $searchRange.Expand($wdParagraph)\n    $paragraphStart = $searchRange.Start\n    if ($paragraphStart -le 0) {\n        return\n    }\n\n    $previousRange = $d.Range($paragraphStart - 1, $paragraphStart - 1)\n    $previousRange.Expand($wdParagraph)\n\n    # A paragraph range that is only the paragraph mark (length 1) is blank.\n    if ($previousRange.Text.Length -eq 1) {\n        $previousRange.Delete()\n    }\n}\n\nRemove-BlankParagraphBefore \"Freudentage, Regenmomente und alles dazwischen\"\nRemove-BlankParagraphBefore \"Gott verhei\u00dft uns seinen Schutz\"\nRemove-BlankParagraphBefore \"Ich finde dieses Motiv des Regenbogens\"\n"}
